# Generate Report for Handoff
# Refresh the "handoff" timestamps recorded for the rows whose status is
# "Handback transform failed" / "Ready for handoff" (i.e. the files that were
# (re)handed-off during this run). Each per-language detail sheet keeps its
# own "Latest Handoff Datetime" and the Overview sheet keeps the rolled-up
# "Latest Handoff Date" - both need to move to the new run time.

$wb = $excel.ActiveWorkbook

$rows = @(7, 10, 11, 12, 13, 14, 15, 16)

# zh-cn sheet: column E is "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "2016-03-19 20:24:05"
}

# de-de sheet: column E is "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "2016-03-19 20:24:11"
}

# Overview sheet: column D is "Latest Handoff Date"
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("D$r").Value = "2016-24-19 20:24:11"
}
